$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.968.75'
$ws.Range("E2").Value = '  +1.58%  '

# Row 3
$ws.Range("D3").Value = '3.858.80'
$ws.Range("E3").Value = '  +1.32%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '473.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.66%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.67%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.85%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
$ws.Range("E9").Value = '  +1.91%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000313'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.48%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.62%  '

# Row 13
$ws.Range("E13").Value = '  -0.59%  '

# Row 14
$ws.Range("D14").Value = '4.479.97'
$ws.Range("E14").Value = '  +1.53%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.37%  '

# Row 16
$ws.Range("D16").Value = '3.856.45'
$ws.Range("E16").Value = '  +0.94%  '

# Row 17
$ws.Range("E17").Value = '  -0.35%  '

# Row 18
$ws.Range("E18").Value = '  +0.28%  '

# Row 19
$ws.Range("E19").Value = '  +3.72%  '

# Row 20
$ws.Range("D20").Value = '67.270.71'
$ws.Range("E20").Value = '  +1.67%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.62%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.82%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.31%  '

# Row 25
$ws.Range("E25").Value = '  +8.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.61%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.31%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.65%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '725.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.05%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.93%  '

# Row 32
$ws.Range("E32").Value = '  +6.19%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.93%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '43.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.67%  '

# Row 35
$ws.Range("E35").Value = '  +8.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.80%  '

# Row 37
$ws.Range("E37").Value = '  +0.11%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.50'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.80%  '

# Row 39
$ws.Range("E39").Value = '  +3.10%  '

# Row 40
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.86%  '

# Row 41
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.346'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.29%  '

# Row 42
$ws.Range("E42").Value = '  +3.23%  '

# Row 43
$ws.Range("D43").Value = '0.0₃0678'
$ws.Range("E43").Value = '  -6.22%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.54'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.79%  '

# Row 46
$ws.Range("E46").Value = '  +2.19%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.69%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.18'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.80%  '

# Row 50
$ws.Range("E50").Value = '  +2.82%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.36%  '
